# Communication menu on dedicated
# 1) Add a new Todo/History row documenting the dedicated-server comms issue.
# 2) Add a "Suggestions" worksheet that lists ideas/modules for the mission.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Todo & History" (was "Tabelle1") - append the new log entry
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Todo & History"

$ws1.Range("A27").Value = 41965
$ws1.Range("A27").NumberFormat = "mm-dd-yy"
$ws1.Range("A27").HorizontalAlignment = -4108

$ws1.Range("B27").Value = "20.45"
$ws1.Range("C27").Value = "On Dedicated, communication menus didn''t appear, needs new testing"

$ws1.Range("A28").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Suggestions" - new sheet right after "Todo & History"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Suggestions"

$ws2.Columns("A:A").ColumnWidth = 28.5703125
$ws2.Columns("B:B").ColumnWidth = 190.140625

# Header row
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Description"

# Column A (Name) values, entered top to bottom
$ws2.Range("A2").Value = "Modules"
$ws2.Range("A3").Value = "Zombie Module"
$ws2.Range("A4").Value = "Patrolling Enemies"
$ws2.Range("A5").Value = "Apocalypse Module"
$ws2.Range("A6").Value = "More Advanced Start Options"
$ws2.Range("A7").Value = "Mod Content as Modules"

# Column B (Description) values for rows 3-6 first ...
$ws2.Range("B3").Value = "There exists a zombie mission, at the moment abandoned, which doesn''t require a mod to spawn (low quality) zombies. Zombies could just spawn around houses or around the player and be an extra threat."
$ws2.Range("B4").Value = "Was planed but never implemented. Units can be created at random which will patrol from one target to another. Important is that they don''t patrol around conquered targets, so that the players are safe around them. We might make an option to disable AT and AA soldiers, so it''s still somewhat safe to drive around."
$ws2.Range("B5").Value = "I heard about a mod which adds apocalyptic stuff to Arma. Apparently there are catastrophes like tornados or earthquakes, which could be another threat."
$ws2.Range("B6").Value = "Additional to the existing advanced start, there should be more options. For example that one assault rifle or launcher is granted or just a vehicle without equipment."

# ... then B2 and B7 last (matches authoring order of the source workbook)
$ws2.Range("B2").Value = "For more variety and customization to the mission, the host should be able to act- resp. deactivate content as modules. For each module should exist a parameter."
$ws2.Range("B7").Value = "Mod Content (like Weapons and Vehicles) can also be activated by parameters. There should be a parameter for each mod."

# Formatting: rows 2-21 are tall, centered/middle aligned, column B wraps text
$ws2.Range("A2:A21").HorizontalAlignment = -4108
$ws2.Range("A2:A21").VerticalAlignment = -4108
$ws2.Range("A2:B21").RowHeight = 30
$ws2.Range("B2:B21").VerticalAlignment = -4108
$ws2.Range("B2:B21").WrapText = $true

$ws2.Range("A22:A23").RowHeight = 27.75
$ws2.Range("A24:A34").RowHeight = 24.75

$ws2.Range("A1").Select()
